$d = $word.ActiveDocument

# 1) "Tuần 7" -> "Tuần 8" in the week-header paragraph
$d.Content.Find.Execute("Tuần 7", $true, $false, $false, $false, $false, $true, 1, $false, "Tuần 8", 2) | Out-Null

# 2) Append four new paragraphs at the very end of the document:
#    "Tuần 9:" (bold header)
#    "- Bổ sung các chức năng trong database"
#    "+ Liên kết các bảng lại với nhau"
#    "+ Đổ dữ liệu vào trang home"
#
# New blank paragraphs are created first (while formatting is still the
# plain, non-bold style inherited from the last existing paragraph), and
# only afterwards is the header paragraph turned bold - this keeps the
# three body paragraphs free of any bold/bCs run properties, matching how
# Word itself would record the edit.

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$rng = $lastPara.Range.Duplicate
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$paraB_idx = $d.Paragraphs.Count

$rB = $d.Paragraphs.Item($paraB_idx).Range.Duplicate
$rB.Collapse(-1)
$rB.InsertParagraphBefore()
$paraHeader_idx = $paraB_idx
$paraB_idx = $paraB_idx + 1

$rBend = $d.Paragraphs.Item($paraB_idx).Range.Duplicate
$rBend.Collapse(0)
$rBend.InsertParagraphAfter()
$paraC_idx = $d.Paragraphs.Count

$rCend = $d.Paragraphs.Item($paraC_idx).Range.Duplicate
$rCend.Collapse(0)
$rCend.InsertParagraphAfter()
$paraD_idx = $d.Paragraphs.Count

$d.Paragraphs.Item($paraHeader_idx).Range.InsertBefore("Tuần 9:")
$d.Paragraphs.Item($paraB_idx).Range.InsertBefore("- Bổ sung các chức năng trong database")
$d.Paragraphs.Item($paraC_idx).Range.InsertBefore("+ Liên kết các bảng lại với nhau")
$d.Paragraphs.Item($paraD_idx).Range.InsertBefore("+ Đổ dữ liệu vào trang home")

$hRng = $d.Paragraphs.Item($paraHeader_idx).Range
$hRng.Bold = 1
$hRng.BoldBi = 1
